$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Fix style of a handful of pre-existing cells that used the now-removed
#     style (fontId 2 / empty font). They should use style index 1 instead,
#     same as their neighboring cells in the same rows.
$ws.Range("C168").Copy()
$ws.Range("B168").PasteSpecial($xlPasteFormats)

$ws.Range("O169").Copy()
$ws.Range("P169:Q169").PasteSpecial($xlPasteFormats)

$ws.Range("O171").Copy()
$ws.Range("P171:Q171").PasteSpecial($xlPasteFormats)

# --- Append new rows of data (172-175), re-using the formatting of
#     existing rows so no new style/font entries are introduced.
$ws.Range("A163:Q163").Copy()
$ws.Range("A172:Q172").PasteSpecial($xlPasteFormats)

$ws.Range("A163:Q163").Copy()
$ws.Range("A173:Q173").PasteSpecial($xlPasteFormats)

$ws.Range("A163:Q163").Copy()
$ws.Range("A174:Q174").PasteSpecial($xlPasteFormats)

$ws.Range("A166:Q166").Copy()
$ws.Range("A175:Q175").PasteSpecial($xlPasteFormats)

# Row 172: 2020-10-18
$ws.Range("A172").Value = 44122
$ws.Range("B172").Value = 847926
$ws.Range("C172").Value = 133841
$ws.Range("D172").Value = 713802
$ws.Range("E172").Value = 283
$ws.Range("F172").Value = 2175
$ws.Range("G172").Value = 125601
$ws.Range("H172").Value = 232
$ws.Range("I172").Value = 76
$ws.Range("J172").Value = 46
$ws.Range("K172").Value = 0
$ws.Range("L172").Value = 4
$ws.Range("M172").Value = 828
$ws.Range("N172").Value = 1176
$ws.Range("O172").Value = 38
$ws.Range("P172").Value = 123
$ws.Range("Q172").Value = 10

# Row 173: 2020-10-19
$ws.Range("A173").Value = 44123
$ws.Range("B173").Value = 849137
$ws.Range("C173").Value = 133932
$ws.Range("D173").Value = 715008
$ws.Range("E173").Value = 197
$ws.Range("F173").Value = 2178
$ws.Range("G173").Value = 126075
$ws.Range("H173").Value = 221
$ws.Range("I173").Value = 70
$ws.Range("J173").Value = 47
$ws.Range("K173").Value = 0
$ws.Range("L173").Value = 3
$ws.Range("M173").Value = 830
$ws.Range("N173").Value = 1177
$ws.Range("O173").Value = 38
$ws.Range("P173").Value = 123
$ws.Range("Q173").Value = 10

# Row 174: 2020-10-20
$ws.Range("A174").Value = 44124
$ws.Range("B174").Value = 852760
$ws.Range("C174").Value = 134324
$ws.Range("D174").Value = 718241
$ws.Range("E174").Value = 195
$ws.Range("F174").Value = 2180
$ws.Range("G174").Value = 126734
$ws.Range("H174").Value = 221
$ws.Range("I174").Value = 86
$ws.Range("J174").Value = 47
$ws.Range("K174").Value = 0
$ws.Range("L174").Value = 2
$ws.Range("M174").Value = 832
$ws.Range("N174").Value = 1177
$ws.Range("O174").Value = 38
$ws.Range("P174").Value = 123
$ws.Range("Q174").Value = 10

# Row 175: 2020-10-21
$ws.Range("A175").Value = 44125
$ws.Range("B175").Value = 856670
$ws.Range("C175").Value = 134588
$ws.Range("D175").Value = 721954
$ws.Range("E175").Value = 128
$ws.Range("F175").Value = 2184
$ws.Range("G175").Value = 127120
$ws.Range("H175").Value = 226
$ws.Range("I175").Value = 86
$ws.Range("J175").Value = 47
$ws.Range("K175").Value = 0
$ws.Range("L175").Value = 2
$ws.Range("M175").Value = 834
$ws.Range("N175").Value = 1179
$ws.Range("O175").Value = 38
$ws.Range("P175").Value = 123
$ws.Range("Q175").Value = 10

Write-Output "done"
